# Textbox response formatting fix
# Rename task-order sheets (new timestamped run) and update the stimulus
# file names / response labels listed on each sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG_TO-... ---------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-1651168669926379"
$ws1.Range("B2").Value = "go_stims-16511686698866572.csv"
$ws1.Range("B3").Value = "GNG_stims-16511686699108255.csv"
$ws1.Range("B4").Value = "go_stims-16511686699118295.csv"
$ws1.Range("B5").Value = "GNG_stims-16511686699253788.csv"

# --- Sheet 2: NB_TO-... -----------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16511686740521245"
$ws2.Range("B2").Value = "ZB-match_8-16511686701473618.csv"
$ws2.Range("B3").Value = "TB-16511686740250905.csv"
$ws2.Range("B4").Value = "ZB-match_1-1651168670245801.csv"
$ws2.Range("B5").Value = "OB-16511686711908796.csv"
$ws2.Range("B6").Value = "ZB-match_7-16511686699808261.csv"
$ws2.Range("B7").Value = "TB-16511686717115793.csv"
$ws2.Range("B8").Value = "OB-16511686709281902.csv"
$ws2.Range("B9").Value = "TB-16511686735301046.csv"
$ws2.Range("B10").Value = "OB-1651168670329121.csv"

# --- Sheet 3: RS_TO-... -------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16511686740521245"
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# --- Sheet 4: TOL_TO-... -------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16511686741009889"
$ws4.Range("B2").Value = "MM_stims-16511686740685406.csv"
$ws4.Range("B3").Value = "ZM_stims-1651168674055434.csv"
$ws4.Range("B4").Value = "MM_stims-16511686740842102.csv"
$ws4.Range("B5").Value = "ZM_stims-16511686740695448.csv"
$ws4.Range("B6").Value = "MM_stims-16511686741000202.csv"
$ws4.Range("B7").Value = "ZM_stims-16511686740851789.csv"

# --- Sheet 5: vSAT_TO-... -------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16511686741615906"
$ws5.Range("B2").Value = "SAT_stims-16511686741040103.csv"
$ws5.Range("B3").Value = "vSAT_stims-16511686741468241.csv"
$ws5.Range("B4").Value = "vSAT_stims-16511686741317866.csv"
$ws5.Range("B5").Value = "SAT_stims-16511686741154182.csv"
